# WS_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer (A16)
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-13
#
# The sheet ships protected, so unlock it for the edit and restore
# protection afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Confidential disclosure banner: refresh the "as of" date ---
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) refresh, rows 2-13 ---
$values = @{
    2  = @(0.02752060754815729, -0.01484523057485787)
    3  = @(0.02134641475285731, -0.009043312708234219)
    4  = @(0.05626675149991613, -0.002112676056338025)
    5  = @(0.1390119078878892,  -0.003703132232680151)
    6  = @(0.02024167090193203, -0.008990318118948903)
    7  = @(0.1294602650075715,  -0.01203369434416368)
    8  = @(0.08887261262550855, -0.01024327784891166)
    9  = @(0.0295167327036918,  -0.01220818162347403)
    10 = @(0.103832375158639,   -0.01305334846765038)
    11 = @(0.2945541072741727,   0.002162341982701266)
    12 = @(0.08937655463966453, -0.006002182611858919)
    13 = @(1,                   -0.005500689646452339)
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# --- Restore the sheet protection that was in place before the edit ---
$ws.Protect()
